$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D12").Value = "“파이토치로 배우는 자연어 처리”가 출간되었습니다."
$ws.Range("E12").Value = "https://tensorflow.blog/2021/05/26/%ed%8c%8c%ec%9d%b4%ed%86%a0%ec%b9%98%eb%a1%9c-%eb%b0%b0%ec%9a%b0%eb%8a%94-%ec%9e%90%ec%97%b0%ec%96%b4-%ec%b2%98%eb%a6%ac%ea%b0%80-%ec%b6%9c%ea%b0%84%eb%90%98%ec%97%88%ec%8a%b5%eb%8b%88%eb%8b%a4/"

$ws.Range("D21").Value = "[Kaldi] Voxceleb Recipe로 i-vector 와 x-vector 성능 비교"
$ws.Range("E21").Value = "https://ms-review.tistory.com/13"

$ws.Range("D32").Value = "KS test (Kolmogorov–Smirnov test)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/320"

$ws.Range("D37").Value = "[Rehearsal] 2021 대한산업공학회 춘계 학술대회 - 김정희"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=1542&mod=document&pageid=1"

$ws.Range("D51").Value = "[css] 자간 조절하기, letter-spacing 속성"
$ws.Range("E51").Value = "https://bskyvision.com/1202"
